$d = $word.ActiveDocument

# Replace the unique sequence "BC Stats BC Data Services" (spanning three
# runs: "BC Stats", " ", "BC Data Services") with a single combined run
# "BC Stats - BC Data Services". Using the full unique phrase avoids
# matching the other standalone "BC Stats" occurrences elsewhere in the
# document.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("BC Stats BC Data Services", $true, $false, $false, $false, $false, `
               $true, 1, $false, "BC Stats - BC Data Services", 2)

$d.Save()
